# Auto-generated Excel COM-interop script to apply F-column (想去人数) updates
# across sheets 展览, 演出, 本地生活, 全部类型 per the target diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 972  # was 970 - 上海·Look Look动漫嘉年华
$ws.Cells.Item(8, 6).Value = 957  # was 953 - 上海·第五人格ONLY
$ws.Cells.Item(11, 6).Value = 1054  # was 1053 - 上海·Redamancy动漫游戏嘉年华-相遇篇
$ws.Cells.Item(12, 6).Value = 4028  # was 3273 - 上海·城市动漫节
$ws.Cells.Item(15, 6).Value = 1705  # was 1703 - 上海·第十九届次元漫步动漫游戏嘉年华
$ws.Cells.Item(16, 6).Value = 23  # was 22 - 上海·第十届ACBC动漫盛典
$ws.Cells.Item(17, 6).Value = 652  # was 649 - 上海·蔚蓝档案ONLY01
$ws.Cells.Item(18, 6).Value = 23  # was 21 - 上海·风花节-花与云之诗
$ws.Cells.Item(20, 6).Value = 379  # was 378 - 上海·第二届奇卡波利国潮嘉年华
$ws.Cells.Item(21, 6).Value = 1095  # was 1094 - 上海·第五十四届妖漫动漫游戏展
$ws.Cells.Item(22, 6).Value = 1531  # was 1528 - 上海·首届Redamancy动漫游戏嘉年华
$ws.Cells.Item(23, 6).Value = 775  # was 773 - 上海·原神X星穹铁道ONLY
$ws.Cells.Item(24, 6).Value = 665  # was 660 - 上海·首届sunshine跨次元动漫游戏展
$ws.Cells.Item(25, 6).Value = 519  # was 517 - 上海·魔都coser动漫展-C展
$ws.Cells.Item(28, 6).Value = 52  # was 49 - 上海·重返未来1999ONLY·2024—UTTU闪烁集会「春申魅影」
$ws.Cells.Item(29, 6).Value = 1034  # was 1031 - 上海·第三届奇卡波利国潮嘉年华
$ws.Cells.Item(30, 6).Value = 1168  # was 1166 - 上海·第二届Redamancy动漫游戏嘉年华
$ws.Cells.Item(31, 6).Value = 339  # was 336 - 上海·运动番ONLY
$ws.Cells.Item(32, 6).Value = 2467  # was 2462 - 上海·Virtual Shanghai Anime Exhibition魔都虚幻世界二次元1.0
$ws.Cells.Item(33, 6).Value = 284  # was 283 - 上海·S·CGE动漫游戏嘉年华
$ws.Cells.Item(34, 6).Value = 1437  # was 1423 - 上海·第三届Redamancy动漫游戏嘉年华
$ws.Cells.Item(36, 6).Value = 8  # was 6 - 上海·第五十八届燃梦星辰国潮嘉年华-随机宅舞
$ws.Cells.Item(37, 6).Value = 68  # was 66 - 上海·灌篮高手--青春永不散场
$ws.Cells.Item(38, 6).Value = 4057  # was 4049 - 上海·原神×崩坏×星铁only旅行盛宴2.0
$ws.Cells.Item(39, 6).Value = 61  # was 59 - 上海·第八届ACBC动漫盛典-国潮汉服游园会

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 111  # was 110 - 上海·三月的幻想演唱会2024「飞越蓝色时刻」
$ws.Cells.Item(6, 6).Value = 198  # was 197 - 上海·遇见新海诚--帝玖「这次一定」室内乐ACG音乐会
$ws.Cells.Item(8, 6).Value = 13  # was 11 - 上海·「再现经典」古典乐巨匠之夜——贝多芬传世经典音乐会《命运交响曲》
$ws.Cells.Item(9, 6).Value = 10  # was 8 - 上海·《卡农Canon in D》世界经典作品视听音乐会
$ws.Cells.Item(14, 6).Value = 4144  # was 4142 - 上海·RAISE A SUILEN ASIA TOUR 2024 IN SHANGHAI
$ws.Cells.Item(15, 6).Value = 6  # was 0 - 上海·「极致现场」草原天籁之声——NAIR奈热乐队音乐会
$ws.Cells.Item(16, 6).Value = 21  # was 19 - 上海·《天空之城》宫崎骏&久石让经典作品动漫视听音乐会
$ws.Cells.Item(23, 6).Value = 264  # was 263 - 上海·《四月是你的谎言》——“公生”与“薰”的钢琴小提琴唯美经典音乐集
$ws.Cells.Item(25, 6).Value = 130  # was 129 - 上海· 茅原实里动漫交响音乐会
$ws.Cells.Item(35, 6).Value = 2  # was 0 - 上海·【520矩献】《爱乐之城》唯美浪漫经典爱情影视听音乐会
$ws.Cells.Item(37, 6).Value = 18  # was 17 - 上海·菊次郎的夏天——久石让钢琴曲梦幻之旅演奏会
$ws.Cells.Item(38, 6).Value = 14  # was 10 - 上海·「多厨狂喜」白金交响乐团二次元交响音乐会

# ---- Sheet: 本地生活 ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 51  # was 50 - 上海·方块大战（豫园店）
$ws.Cells.Item(4, 6).Value = 1293  # was 1290 - 上海·罗小黑 x HAPPY ZOO主题Cafe
$ws.Cells.Item(5, 6).Value = 1689  # was 1687 - 上海・明日方舟主题店·[SWEET ZONE甜蜜区域]
$ws.Cells.Item(8, 6).Value = 89  # was 65 - 上海·NIJISANJI EN 官方授权主题店

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 1293  # was 1290 - 上海·罗小黑 x HAPPY ZOO主题Cafe
$ws.Cells.Item(3, 6).Value = 1689  # was 1687 - 上海・明日方舟主题店·[SWEET ZONE甜蜜区域]
$ws.Cells.Item(7, 6).Value = 972  # was 970 - 上海·Look Look动漫嘉年华
$ws.Cells.Item(10, 6).Value = 957  # was 953 - 上海·第五人格ONLY
$ws.Cells.Item(14, 6).Value = 198  # was 197 - 上海·遇见新海诚--帝玖「这次一定」室内乐ACG音乐会
$ws.Cells.Item(15, 6).Value = 198  # was 197 - 上海·遇见新海诚--帝玖「这次一定」室内乐ACG音乐会
$ws.Cells.Item(16, 6).Value = 13  # was 11 - 上海·「再现经典」古典乐巨匠之夜——贝多芬传世经典音乐会《命运交响曲》
$ws.Cells.Item(17, 6).Value = 1054  # was 1053 - 上海·Redamancy动漫游戏嘉年华-相遇篇
$ws.Cells.Item(18, 6).Value = 10  # was 8 - 上海·《卡农Canon in D》世界经典作品视听音乐会
$ws.Cells.Item(19, 6).Value = 4028  # was 3299 - 上海·城市动漫节
$ws.Cells.Item(22, 6).Value = 1705  # was 1703 - 上海·第十九届次元漫步动漫游戏嘉年华
$ws.Cells.Item(23, 6).Value = 23  # was 22 - 上海·第十届ACBC动漫盛典
$ws.Cells.Item(24, 6).Value = 652  # was 649 - 上海·蔚蓝档案ONLY01
$ws.Cells.Item(26, 6).Value = 379  # was 378 - 上海·第二届奇卡波利国潮嘉年华
$ws.Cells.Item(27, 6).Value = 1095  # was 1094 - 上海·第五十四届妖漫动漫游戏展
$ws.Cells.Item(28, 6).Value = 1531  # was 1528 - 上海·首届Redamancy动漫游戏嘉年华
$ws.Cells.Item(31, 6).Value = 775  # was 773 - 上海·原神X星穹铁道ONLY
$ws.Cells.Item(32, 6).Value = 665  # was 660 - 上海·首届sunshine跨次元动漫游戏展
$ws.Cells.Item(33, 6).Value = 519  # was 517 - 上海·魔都coser动漫展-C展
$ws.Cells.Item(36, 6).Value = 52  # was 49 - 上海·重返未来1999ONLY·2024—UTTU闪烁集会「春申魅影」
$ws.Cells.Item(39, 6).Value = 264  # was 263 - 上海·《四月是你的谎言》——“公生”与“薰”的钢琴小提琴唯美经典音乐集
$ws.Cells.Item(40, 6).Value = 1034  # was 1031 - 上海·第三届奇卡波利国潮嘉年华
$ws.Cells.Item(41, 6).Value = 1168  # was 1166 - 上海·第二届Redamancy动漫游戏嘉年华
$ws.Cells.Item(42, 6).Value = 339  # was 336 - 上海·运动番ONLY
$ws.Cells.Item(43, 6).Value = 2467  # was 2462 - 上海·Virtual Shanghai Anime Exhibition魔都虚幻世界二次元1.0
$ws.Cells.Item(46, 6).Value = 1437  # was 1423 - 上海·第三届Redamancy动漫游戏嘉年华
$ws.Cells.Item(48, 6).Value = 8  # was 6 - 上海·第五十八届燃梦星辰国潮嘉年华-随机宅舞
$ws.Cells.Item(50, 6).Value = 4057  # was 4049 - 上海·原神×崩坏×星铁only旅行盛宴2.0
$ws.Cells.Item(51, 6).Value = 14  # was 10 - 上海·「多厨狂喜」白金交响乐团二次元交响音乐会
